$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.523.68'
$ws.Range('E2').Value = '  -5.72%  '
$ws.Range('D3').Value = '2.980.66'
$ws.Range('E3').Value = '  -7.12%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '542.52'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -5.72%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '152.50'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -9.29%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.560'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -6.40%  '
$ws.Range('D9').Value = '2.987.36'
$ws.Range('E9').Value = '  -6.77%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.111'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -6.87%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.18'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -8.53%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.364'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -7.52%  '
$ws.Range('D13').Value = '3.504.50'
$ws.Range('E13').Value = '  -7.11%  '
$ws.Range('E14').Value = '  -3.73%  '
$ws.Range('D15').Value = '61.656.11'
$ws.Range('E15').Value = '  -5.45%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '23.52'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -8.60%  '
$ws.Range('D17').Value = '2.984.57'
$ws.Range('E17').Value = '  -6.72%  '
$ws.Range('E18').Value = '  -7.31%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.10'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.76%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '384.21'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -7.24%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.89'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -8.21%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.60'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -8.59%  '
$ws.Range('E23').Value = '  -0.10%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.92'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -6.71%  '
$ws.Range('E25').Value = '  -5.49%  '
$ws.Range('E26').Value = '  -7.87%  '
$ws.Range('E27').Value = '  +0.17%  '
$ws.Range('D28').Value = '0.0₃0931'
$ws.Range('E28').Value = '  -11.89%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.26'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -7.58%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.999'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.01%  '
$ws.Range('E31').Value = '  -7.70%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.21'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -6.54%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '158.45'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.23%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.99'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -6.98%  '
$ws.Range('E35').Value = '  -8.05%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.06'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -7.63%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.27'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -7.52%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.57'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -8.88%  '
$ws.Range('D39').Value = '2.429.05'
$ws.Range('E39').Value = '  -11.85%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.87'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -7.15%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '37.05'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -5.23%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '22.16'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -9.22%  '
$ws.Range('E43').Value = '  -7.37%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0591'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -6.83%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.999'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0244'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -7.71%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.95'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -13.42%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0953'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.30%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '19.62'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -9.70%  '
$ws.Range('E50').Value = '  +0.26%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '262.93'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -11.76%  '
